$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Iteration 1 block (rows 6-10): update names/comments ---
# Row 9: Alfred -> Xuxin
$ws.Range("B9").Value = "Xuxin"

# Comments (column D) for rows 6-9 change text
$ws.Range("D6").Value = "Styling for parties list and adding testing"
$ws.Range("D7").Value = "Worked on reviewing and fitting git pull requests"
$ws.Range("D8").Value = "Potential drag and drop feature for voting"
$ws.Range("D9").Value = "Created voting method for above and low the line"

# --- Iteration 2 block (rows 14-18): fill in grade + comment columns ---
$ws.Range("C14").Value = "D"
$ws.Range("D14").Value = "Styling for parties list and adding testing"

$ws.Range("C15").Value = "HD"
$ws.Range("D15").Value = "Worked on reviewing and fitting git pull requests"

$ws.Range("C16").Value = "D"
$ws.Range("D16").Value = "Potential drag and drop feature for voting"

$ws.Range("B17").Formula = "=B9"
$ws.Range("C17").Value = "HD"
$ws.Range("D17").Value = "Created voting method for above and low the line"

$ws.Range("C18").Value = "C"
$ws.Range("D18").Value = "Assisting with basic functionality "

# --- Update view: active cell selection (scrolled so row 10 is visible) ---
$excel.Goto($ws.Range("A10"), $false)
$ws.Range("B18").Select()
